$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Add the three new data rows for years 2001-2003 (A/B/C/D columns)
$ws.Cells.Item(53,1).Value = 404
$ws.Cells.Item(53,2).Value = "Kenya"
$ws.Cells.Item(53,3).Value = "Historical Gender Equality Index"
$ws.Cells.Item(53,4).Value = 2001

$ws.Cells.Item(54,1).Value = 404
$ws.Cells.Item(54,2).Value = "Kenya"
$ws.Cells.Item(54,3).Value = "Historical Gender Equality Index"
$ws.Cells.Item(54,4).Value = 2002

$ws.Cells.Item(55,1).Value = 404
$ws.Cells.Item(55,2).Value = "Kenya"
$ws.Cells.Item(55,3).Value = "Historical Gender Equality Index"
$ws.Cells.Item(55,4).Value = 2003

# Update the Data column (E) for rows 2-55 with the new Historical Gender Equality
# Index values. Each is written as a quoted-string formula first, then pasted back
# as a value-only paste so Excel stores it as shared-string text (matching the
# source data, which stores these look-like-numbers as text) without leaving any
# numeric-format/style residue behind.
$ws.Cells.Item(2,5).Formula = "=""58.337616656494"""
$ws.Cells.Item(3,5).Formula = "=""59.8020731992684"""
$ws.Cells.Item(4,5).Formula = "=""59.7931203068408"""
$ws.Cells.Item(5,5).Formula = "=""60.7068799963169"""
$ws.Cells.Item(6,5).Formula = "=""61.1706847624426"""
$ws.Cells.Item(7,5).Formula = "=""57.7230423970673"""
$ws.Cells.Item(8,5).Formula = "=""59.6225548866545"""
$ws.Cells.Item(9,5).Formula = "=""60.5908203647597"""
$ws.Cells.Item(10,5).Formula = "=""59.3392360562427"""
$ws.Cells.Item(11,5).Formula = "=""59.9470353293592"""
$ws.Cells.Item(12,5).Formula = "=""56.9969886684851"""
$ws.Cells.Item(13,5).Formula = "=""60.6127075897675"""
$ws.Cells.Item(14,5).Formula = "=""59.4789541171836"""
$ws.Cells.Item(15,5).Formula = "=""58.9865926112627"""
$ws.Cells.Item(16,5).Formula = "=""61.7170264797121"""
$ws.Cells.Item(17,5).Formula = "=""57.4533690447227"""
$ws.Cells.Item(18,5).Formula = "=""58.42642778701"""
$ws.Cells.Item(19,5).Formula = "=""58.308518203806"""
$ws.Cells.Item(20,5).Formula = "=""58.3757491925908"""
$ws.Cells.Item(21,5).Formula = "=""57.148667752659"""
$ws.Cells.Item(22,5).Formula = "=""58.7184948650439"""
$ws.Cells.Item(23,5).Formula = "=""60.1832179418252"""
$ws.Cells.Item(24,5).Formula = "=""57.3143708170811"""
$ws.Cells.Item(25,5).Formula = "=""60.8037889303445"""
$ws.Cells.Item(26,5).Formula = "=""62.0714610907295"""
$ws.Cells.Item(27,5).Formula = "=""59.1444054107502"""
$ws.Cells.Item(28,5).Formula = "=""58.7335804567874"""
$ws.Cells.Item(29,5).Formula = "=""60.5578439662464"""
$ws.Cells.Item(30,5).Formula = "=""61.5737152314073"""
$ws.Cells.Item(31,5).Formula = "=""60.7015897017626"""
$ws.Cells.Item(32,5).Formula = "=""61.427024692804"""
$ws.Cells.Item(33,5).Formula = "=""60.5190980300702"""
$ws.Cells.Item(34,5).Formula = "=""60.7711140979584"""
$ws.Cells.Item(35,5).Formula = "=""58.7127513804857"""
$ws.Cells.Item(36,5).Formula = "=""60.4576735411756"""
$ws.Cells.Item(37,5).Formula = "=""62.3357596926875"""
$ws.Cells.Item(38,5).Formula = "=""60.0995453214294"""
$ws.Cells.Item(39,5).Formula = "=""61.6737902447174"""
$ws.Cells.Item(40,5).Formula = "=""61.032474501975"""
$ws.Cells.Item(41,5).Formula = "=""65.7797571423626"""
$ws.Cells.Item(42,5).Formula = "=""66.2656786762754"""
$ws.Cells.Item(43,5).Formula = "=""63.1340175795683"""
$ws.Cells.Item(44,5).Formula = "=""64.0394003431104"""
$ws.Cells.Item(45,5).Formula = "=""62.0437687905792"""
$ws.Cells.Item(46,5).Formula = "=""62.9122995212335"""
$ws.Cells.Item(47,5).Formula = "=""64.6048415531091"""
$ws.Cells.Item(48,5).Formula = "=""61.0317050095811"""
$ws.Cells.Item(49,5).Formula = "=""63.347165233731"""
$ws.Cells.Item(50,5).Formula = "=""65.0303494849141"""
$ws.Cells.Item(51,5).Formula = "=""62.158042334143"""
$ws.Cells.Item(52,5).Formula = "=""64.2012709238978"""
$ws.Cells.Item(53,5).Formula = "=""63.3036754591565"""
$ws.Cells.Item(54,5).Formula = "=""68.1713461928154"""
$ws.Cells.Item(55,5).Formula = "=""66.4282881860486"""

$ws.Range("E2:E55").Copy()
$ws.Range("E2:E55").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Fix punctuation in the citation text on the Metadata sheet (Oxford comma added
# before "and Auke Rijpma").
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(3,2).Value = "Carmichael, Sarah, Selin Dilli, and Auke Rijpma (2015). Historical Gender Equality Index. http://hdl.handle.net/10622/VHYIAT, accessed via the Clio Infra website."
